$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = @(0.2917716402565462, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 1, 4.630591052383734)
    3 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    4 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 8.974608811992548)
    5 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    6 = @(0.0006408296065709695, 0.04071648406533734, 0.1494219747398047, 0.4942365360607697, 1, 0.6850158244724827)
    7 = @(0.2917716402565462, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 1, 5.979547907384307)
}

foreach ($row in $values.Keys) {
    $rowValues = $values[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowValues[$col - 2]
    }
}
